# update template because box file path changed

$wb = $excel.ActiveWorkbook

# --- Rename the "text" sheet to "path" ---
$wsPath = $wb.Worksheets.Item("text")
$wsPath.Name = "path"

# --- Make "path" the active/selected sheet (it was "model" before) ---
$wsPath.Activate()

# --- Update the Box-backed file paths on the "path" sheet ---
# The Box desktop app moved from "~/Box/..." to
# "/Users/michaelfive/Box/Box 3EA Team Folder/..." paths.
# C2 = model_file_path (Niger PSRA folder) is written first, then
# B2 = data_file_path (Lebanon preimputed data file), so the new text
# lands in shared-strings in that same order.
$wsPath.Range("C2").Value = "/Users/michaelfive/Box/Box 3EA Team Folder/For Zezhen/MR automation/Test Data_Niger/PSRA"
$wsPath.Range("B2").Value = "/Users/michaelfive/Box/Box 3EA Team Folder/3EA Analysis/3EA Lebanon_Analysis/Lebanon_Y1_FA/LBY1_PREIMPUTED_FULL_SPREAD_10-31-2019_mplus.dta"

# --- Move the selection on the "path" sheet from C3 to C2 ---
$wsPath.Range("C2").Select()
